# Rename the 'Middagstudies' column header to the more generic 'Sancties'
# on both class sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Klas1")
$ws1.Range("B1").Value = "Sancties"

$ws2 = $wb.Worksheets.Item("Klas2")
$ws2.Range("B1").Value = "Sancties"

# Restore the cursor/selection state left behind by the interactive edit.
$ws1.Activate()
$ws1.Columns("C").Select()

$ws2.Activate()
$ws2.Range("G16").Select()
